# Fill marker info ("NAT") into column J for the rows that were
# still missing it, and move the active cell selection to J20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(11, 12, 13, 14, 15, 16, 17, 18, 19, 23, 24, 25)
foreach ($r in $rows) {
    $ws.Range("J$r").Value = "NAT"
}

$ws.Range("J20").Select()
